$d = $word.ActiveDocument

$d.Content.Find.Execute('66÷7=9, 3', $true, $true, $false, $false, $false, $true, 1, $false, '10÷7=1, 3', 2) | Out-Null
$d.Content.Find.Execute('43÷3=14, 1', $true, $true, $false, $false, $false, $true, 1, $false, '83÷7=11, 6', 2) | Out-Null
$d.Content.Find.Execute('55÷8=6, 7', $true, $true, $false, $false, $false, $true, 1, $false, '80÷7=11, 3', 2) | Out-Null
$d.Content.Find.Execute('93÷6=15, 3', $true, $true, $false, $false, $false, $true, 1, $false, '44÷8=5, 4', 2) | Out-Null
$d.Content.Find.Execute('43÷8=5, 3', $true, $true, $false, $false, $false, $true, 1, $false, '48÷3=16, 0', 2) | Out-Null
$d.Content.Find.Execute('11÷6=1, 5', $true, $true, $false, $false, $false, $true, 1, $false, '36÷7=5, 1', 2) | Out-Null
$d.Content.Find.Execute('10÷4=2, 2', $true, $true, $false, $false, $false, $true, 1, $false, '82÷8=10, 2', 2) | Out-Null
$d.Content.Find.Execute('41÷8=5, 1', $true, $true, $false, $false, $false, $true, 1, $false, '89÷6=14, 5', 2) | Out-Null
$d.Content.Find.Execute('76÷7=10, 6', $true, $true, $false, $false, $false, $true, 1, $false, '30÷4=7, 2', 2) | Out-Null
$d.Content.Find.Execute('67÷6=11, 1', $true, $true, $false, $false, $false, $true, 1, $false, '57÷7=8, 1', 2) | Out-Null
$d.Content.Find.Execute('42÷8=5, 2', $true, $true, $false, $false, $false, $true, 1, $false, '78÷9=8, 6', 2) | Out-Null
$d.Content.Find.Execute('79÷8=9, 7', $true, $true, $false, $false, $false, $true, 1, $false, '80÷5=16, 0', 2) | Out-Null
$d.Content.Find.Execute('63÷5=12, 3', $true, $true, $false, $false, $false, $true, 1, $false, '65÷7=9, 2', 2) | Out-Null
$d.Content.Find.Execute('78÷6=13, 0', $true, $true, $false, $false, $false, $true, 1, $false, '54÷6=9, 0', 2) | Out-Null
$d.Content.Find.Execute('73÷9=8, 1', $true, $true, $false, $false, $false, $true, 1, $false, '59÷3=19, 2', 2) | Out-Null
$d.Content.Find.Execute('57÷3=19, 0', $true, $true, $false, $false, $false, $true, 1, $false, '10÷7=1, 3', 2) | Out-Null
$d.Content.Find.Execute('39÷4=9, 3', $true, $true, $false, $false, $false, $true, 1, $false, '81÷8=10, 1', 2) | Out-Null
$d.Content.Find.Execute('61÷8=7, 5', $true, $true, $false, $false, $false, $true, 1, $false, '60÷7=8, 4', 2) | Out-Null
$d.Content.Find.Execute('12÷6=2, 0', $true, $true, $false, $false, $false, $true, 1, $false, '52÷5=10, 2', 2) | Out-Null
$d.Content.Find.Execute('54÷8=6, 6', $true, $true, $false, $false, $false, $true, 1, $false, '98÷9=10, 8', 2) | Out-Null
$d.Content.Find.Execute('68÷8=8, 4', $true, $true, $false, $false, $false, $true, 1, $false, '52÷9=5, 7', 2) | Out-Null
$d.Content.Find.Execute('85÷2=42, 1', $true, $true, $false, $false, $false, $true, 1, $false, '95÷8=11, 7', 2) | Out-Null
$d.Content.Find.Execute('50÷5=10, 0', $true, $true, $false, $false, $false, $true, 1, $false, '51÷5=10, 1', 2) | Out-Null
$d.Content.Find.Execute('36÷8=4, 4', $true, $true, $false, $false, $false, $true, 1, $false, '27÷5=5, 2', 2) | Out-Null
$d.Content.Find.Execute('41÷3=13, 2', $true, $true, $false, $false, $false, $true, 1, $false, '68÷9=7, 5', 2) | Out-Null
